$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Sep 29 11:36:57 EDT 2023"
$ws.Range("B3").Value = "Fri Sep 29 11:37:11 EDT 2023"
$ws.Range("B4").Value = "Fri Sep 29 11:37:24 EDT 2023"
$ws.Range("B5").Value = "Fri Sep 29 11:37:38 EDT 2023"
